$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New device-master rows appended after the existing data (rows 147-156),
# two new Mac-Address batches of 5 devices each (ids 3000166-3000175).

$data = @{
    147 = @{ A=3000166; B="Finger Print Scanner 30"; C="D6-15-AC-80-6B-86"; D="BS563Q2230814"; F=165 }
    148 = @{ A=3000167; B="IRIS Scanner 30";         C="6D-58-E2-DF-74-34"; D="BS563Q2230815"; F=327 }
    149 = @{ A=3000168; B="Web Camera 30";           C="E2-A8-56-86-15-30"; D="BS563Q2230816"; F=736 }
    150 = @{ A=3000169; B="Document Scanner 30";     C="72-E8-B9-FD-63-65"; D="BS563Q2230817"; F=801 }
    151 = @{ A=3000170; B="Printer 30";               C="D3-F3-A4-50-AD-12"; D="BS563Q2230818"; F=920 }
    152 = @{ A=3000171; B="Finger Print Scanner 31"; C="06-16-D0-0B-A6-E4"; D="BS563Q2230819"; F=165 }
    153 = @{ A=3000172; B="IRIS Scanner 31";         C="21-78-45-AC-E9-20"; D="BS563Q2230820"; F=327 }
    154 = @{ A=3000173; B="Web Camera 31";           C="3C-E8-87-99-DB-FA"; D="BS563Q2230821"; F=736 }
    155 = @{ A=3000174; B="Document Scanner 31";     C="BF-55-53-98-40-08"; D="BS563Q2230822"; F=801 }
    156 = @{ A=3000175; B="Printer 31";               C="5A-43-36-46-22-EB"; D="BS563Q2230823"; F=920 }
}

$batch1 = 147..151
$batch2 = 152..156

# Batch 1 was filled in column order Name -> Mac Address -> Serial Number.
foreach ($r in $batch1) { $ws.Cells.Item($r, 2).Value = $data[$r].B }
foreach ($r in $batch1) { $ws.Cells.Item($r, 3).Value = $data[$r].C }
foreach ($r in $batch1) { $ws.Cells.Item($r, 4).Value = $data[$r].D }

# Batch 2 was filled Name -> Serial Number -> Mac Address.
foreach ($r in $batch2) { $ws.Cells.Item($r, 2).Value = $data[$r].B }
foreach ($r in $batch2) { $ws.Cells.Item($r, 4).Value = $data[$r].D }
foreach ($r in $batch2) { $ws.Cells.Item($r, 3).Value = $data[$r].C }

# Remaining columns (id, dspec_id, lang_code, is_active, cr_by, cr_dtimes)
# reuse already-existing shared strings, so fill order doesn't affect the
# shared-string table.
foreach ($r in ($batch1 + $batch2)) {
    $ws.Cells.Item($r, 1).Value = $data[$r].A
    $ws.Cells.Item($r, 6).Value = $data[$r].F
    $ws.Cells.Item($r, 7).Value = "eng"
    $ws.Cells.Item($r, 8).Value = $true
    $ws.Cells.Item($r, 8).HorizontalAlignment = -4131
    $ws.Cells.Item($r, 9).Value = "superadmin"
    $ws.Cells.Item($r, 10).Value = "now()"
}

$ws.Cells.Item(156, 5).Select()
